$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column I: spike_counter header + total spike count formula
$ws.Range("I1").Value = "spike_counter"
$ws.Range("I2").Formula = "=SUM(G2:G152)"

# H2 loses its formula and becomes a plain static 0
$ws.Range("H2").Value = 0

# New E-column leak/integrate-and-fire formulas (V membrane potential)
$ws.Range("E3").Formula = "=IF(H3=0,0,IF(E2>=`$D`$2,0,E2 + `$A`$2*(-1*(E2/`$C`$2) + `$F`$2 * (H3 / `$C`$2))))"
$ws.Range("E4:E67").Formula = "=IF(H4=0,0,IF(E3>=`$D`$2,0,E3 + `$A`$2*(-1*(E3/`$C`$2) + `$F`$2 * (H4 / `$C`$2))))"
$ws.Range("E68:E131").Formula = "=IF(H68=0,0,IF(E67>=`$D`$2,0,E67 + `$A`$2*(-1*(E67/`$C`$2) + `$F`$2 * (H68 / `$C`$2))))"
$ws.Range("E132:E152").Formula = "=IF(H132=0,0,IF(E131>=`$D`$2,0,E131 + `$A`$2*(-1*(E131/`$C`$2) + `$F`$2 * (H132 / `$C`$2))))"

# View state: zoom + active selection
$ws.Select() | Out-Null
$excel.ActiveWindow.Zoom = 160
$ws.Range("C3").Select() | Out-Null

# Move/resize the chart to make room for the new column
$co = $ws.ChartObjects(1)
$co.Left = 573.3123622047244
$co.Top = 0.0
$co.Width = 320.6248818897638
$co.Height = 216.0
